$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: was inline string "85" -> numeric 85
$ws.Range("D2").Value = 85

# E2: "officer" -> "Member"
$ws.Range("E2").Value = "Member"
